# Added the validation of Shopping Car
# Populate the "Global" sheet with the ShCar Quantity helper column (B)
# that mirrors the existing Quantity column (A) with "SPAN" markers and
# a thin top/bottom border, matching the html-tag validation helper data.

$wb = $excel.ActiveWorkbook
$wsGlobal = $wb.Worksheets.Item("Global")
$wsAction = $wb.Worksheets.Item("Action1")

# Header for the new column
$wsGlobal.Range("B1").Value = "ShCar Quantity_html tag"

# Data rows: B2:B4 get the "SPAN" marker with a thin top+bottom border
foreach ($addr in @("B2", "B3", "B4")) {
    $cell = $wsGlobal.Range($addr)
    $cell.Value = "SPAN"

    $top = $cell.Borders.Item(8)
    $top.LineStyle = 1
    $top.Color = 0

    $bottom = $cell.Borders.Item(9)
    $bottom.LineStyle = 1
    $bottom.Color = 0
}

# Widen column B to fit the new header text
$wsGlobal.Columns("B").ColumnWidth = 20.5

# Restore the selection / active cell on the Global sheet without leaving
# it as the active tab (Action1 stays the active sheet, as before).
$wsGlobal.Range("B4").Select()
$wsAction.Activate()
